# Refresh the cryptocurrency price/volume snapshot in columns D (Price) and
# E (Volume(1h)) for rows 2-51, matching the latest scrape results.
# Column D values are written with a leading apostrophe so Excel stores them
# as literal text (preserving formats like trailing zeros "1.00" or the
# thousands-grouped "65.887.91") instead of silently coercing to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.887.91'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '''3.319.31'
$ws.Range("E3").Value = '  +1.30%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''558.16'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '''184.87'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''3.316.21'
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -7.06%  '
$ws.Range("D11").Value = '''0.577'
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").Value = '''45.76'
$ws.Range("E12").Value = '  -3.69%  '
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").Value = '''3.850.52'
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("D15").Value = '''8.43'
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").Value = '''568.17'
$ws.Range("E16").Value = '  -10.75%  '
$ws.Range("D17").Value = '''65.827.17'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").Value = '''3.318.55'
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = '''17.65'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '''10.81'
$ws.Range("E21").Value = '  -4.89%  '
$ws.Range("D22").Value = '''0.889'
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '''18.00'
$ws.Range("E23").Value = '  -2.31%  '
$ws.Range("D24").Value = '''4.97'
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").Value = '''97.51'
$ws.Range("E25").Value = '  -9.57%  '
$ws.Range("D26").Value = '''3.94'
$ws.Range("E26").Value = '  -0.86%  '
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '''9.35'
$ws.Range("E28").Value = '  -2.68%  '
$ws.Range("D29").Value = '''8.47'
$ws.Range("E29").Value = '  -2.86%  '
$ws.Range("D30").Value = '''30.39'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").Value = '''6.69'
$ws.Range("E31").Value = '  +6.69%  '
$ws.Range("E32").Value = '  -10.21%  '
$ws.Range("D33").Value = '''558.66'
$ws.Range("E33").Value = '  +4.53%  '
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("D36").Value = '''3.732.22'
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '''55.55'
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("D39").Value = '''33.62'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("E40").Value = '  -4.23%  '
$ws.Range("D41").Value = '''0.0₃0683'
$ws.Range("E41").Value = '  -6.36%  '
$ws.Range("E42").Value = '  -7.43%  '
$ws.Range("D43").Value = '''2.57'
$ws.Range("E43").Value = '  -6.31%  '
$ws.Range("D44").Value = '''3.32'
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("D47").Value = '''3.00'
$ws.Range("E47").Value = '  -13.39%  '
$ws.Range("D48").Value = '''0.126'
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").Value = '''0.999'
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = '''2.50'
$ws.Range("E50").Value = '  -4.33%  '
$ws.Range("D51").Value = '''125.13'
$ws.Range("E51").Value = '  +2.78%  '
